# Update generated report output ("Update gh-pages to output generated at 456a3b4")
# Applies refreshed "want to go" (F column) counts across sheets, and replaces
# the oldest "本地生活" (Local Life) entry with a newer one (row shift + new last row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - F column (想去人数 / want-to-go count) refresh
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 338
$ws1.Range("F5").Value  = 78
$ws1.Range("F6").Value  = 2497
$ws1.Range("F7").Value  = 52
$ws1.Range("F10").Value = 1514
$ws1.Range("F12").Value = 611
$ws1.Range("F13").Value = 1471
$ws1.Range("F14").Value = 1471
$ws1.Range("F15").Value = 1219
$ws1.Range("F16").Value = 494
$ws1.Range("F17").Value = 3554
$ws1.Range("F19").Value = 3275
$ws1.Range("F20").Value = 733
$ws1.Range("F21").Value = 2120
$ws1.Range("F23").Value = 285
$ws1.Range("F25").Value = 1114
$ws1.Range("F28").Value = 965
$ws1.Range("F29").Value = 945

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - F column refresh
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 85
$ws2.Range("F20").Value = 173

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life) - the oldest entry
# ("光与夜之恋 × 线条小狗 ×爱胖达文化" at old row 3) has expired/dropped out
# of the feed. The remaining rows shift up by one; the trailing (A column)
# index stays positional (0,1,2,3,4), and the EVANGELION row's want-to-go
# count ticked up from 493 to 495 in the refreshed pull.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Rows("3:3").Delete()
$ws3.Range("A3").Value = 2
$ws3.Range("A4").Value = 3
$ws3.Range("A5").Value = 4
$ws3.Range("F5").Value = 495

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types) - F column refresh (aggregated view)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 338
$ws4.Range("F10").Value = 78
$ws4.Range("F11").Value = 495
$ws4.Range("F12").Value = 2497
$ws4.Range("F14").Value = 52
$ws4.Range("F17").Value = 85
$ws4.Range("F22").Value = 1514
$ws4.Range("F25").Value = 1471
$ws4.Range("F26").Value = 1471
$ws4.Range("F29").Value = 1219
$ws4.Range("F30").Value = 494
$ws4.Range("F32").Value = 3554
$ws4.Range("F34").Value = 3275
$ws4.Range("F35").Value = 733
$ws4.Range("F37").Value = 2120
$ws4.Range("F39").Value = 285
$ws4.Range("F40").Value = 1114
$ws4.Range("F44").Value = 173
$ws4.Range("F48").Value = 965
$ws4.Range("F49").Value = 945
